$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "46.466.89"
$ws.Cells.Item(2, 5).Value = "  +0.57%  "
$ws.Cells.Item(3, 4).Value = "2.588.60"
$ws.Cells.Item(3, 5).Value = "  +9.83%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).Value = "'305.08"
$ws.Cells.Item(5, 5).Value = "  +1.07%  "
$ws.Cells.Item(6, 4).Value = "'99.82"
$ws.Cells.Item(6, 5).Value = "  +0.31%  "
$ws.Cells.Item(7, 5).Value = "  +5.08%  "
$ws.Cells.Item(8, 4).Value = "'0.999"
$ws.Cells.Item(8, 5).Value = "  -0.07%  "
$ws.Cells.Item(9, 4).Value = "'0.581"
$ws.Cells.Item(9, 5).Value = "  +13.37%  "
$ws.Cells.Item(10, 4).Value = "'38.53"
$ws.Cells.Item(10, 5).Value = "  +11.86%  "
$ws.Cells.Item(11, 4).Value = "'0.0836"
$ws.Cells.Item(11, 5).Value = "  +4.50%  "
$ws.Cells.Item(12, 4).Value = "'8.14"
$ws.Cells.Item(12, 5).Value = "  +14.17%  "
$ws.Cells.Item(13, 4).Value = "2.976.84"
$ws.Cells.Item(13, 5).Value = "  +9.57%  "
$ws.Cells.Item(14, 5).Value = "  +1.87%  "
$ws.Cells.Item(15, 4).Value = "2.601.56"
$ws.Cells.Item(15, 5).Value = "  +10.56%  "
$ws.Cells.Item(16, 5).Value = "  +11.13%  "
$ws.Cells.Item(17, 4).Value = "'14.85"
$ws.Cells.Item(17, 5).Value = "  +9.16%  "
$ws.Cells.Item(18, 4).Value = "46.569.58"
$ws.Cells.Item(18, 5).Value = "  +1.01%  "
$ws.Cells.Item(19, 4).Value = "'13.35"
$ws.Cells.Item(19, 5).Value = "  +4.80%  "
$ws.Cells.Item(20, 5).Value = "  +4.57%  "
$ws.Cells.Item(21, 5).Value = "  +9.36%  "
$ws.Cells.Item(22, 4).Value = "'71.27"
$ws.Cells.Item(22, 5).Value = "  +5.82%  "
$ws.Cells.Item(23, 4).Value = "'255.58"
$ws.Cells.Item(23, 5).Value = "  +3.53%  "
$ws.Cells.Item(24, 5).Value = "  +5.02%  "
$ws.Cells.Item(25, 4).Value = "'2.18"
$ws.Cells.Item(25, 5).Value = "  +13.56%  "
$ws.Cells.Item(26, 4).Value = "'27.96"
$ws.Cells.Item(26, 5).Value = "  +32.84%  "
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 5).Value = "  +0.01%  "
$ws.Cells.Item(28, 4).Value = "'10.48"
$ws.Cells.Item(28, 5).Value = "  +6.81%  "
$ws.Cells.Item(29, 5).Value = "  +4.53%  "
$ws.Cells.Item(30, 4).Value = "'39.65"
$ws.Cells.Item(30, 5).Value = "  +0.10%  "
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).Value = "'6.11"
$ws.Cells.Item(31, 5).Value = "  +10.47%  "
$ws.Cells.Item(32, 2).Value = "LidoDAOToken"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(32, 4).Value = "'3.72"
$ws.Cells.Item(32, 5).Value = "  -0.48%  "
$ws.Cells.Item(33, 4).Value = "'2.32"
$ws.Cells.Item(33, 5).Value = "  +22.32%  "
$ws.Cells.Item(34, 5).Value = "  +5.12%  "
$ws.Cells.Item(35, 4).Value = "'0.0831"
$ws.Cells.Item(35, 5).Value = "  +7.10%  "
$ws.Cells.Item(36, 4).Value = "'149.85"
$ws.Cells.Item(36, 5).Value = "  +2.59%  "
$ws.Cells.Item(37, 5).Value = "  +4.14%  "
$ws.Cells.Item(38, 5).Value = "  +4.36%  "
$ws.Cells.Item(39, 4).Value = "'4.18"
$ws.Cells.Item(39, 5).Value = "  +5.87%  "
$ws.Cells.Item(40, 4).Value = "'15.72"
$ws.Cells.Item(40, 5).Value = "  +5.02%  "
$ws.Cells.Item(41, 4).Value = "'3.62"
$ws.Cells.Item(41, 5).Value = "  +12.26%  "
$ws.Cells.Item(42, 5).Value = "  +6.98%  "
$ws.Cells.Item(43, 4).Value = "2.025.78"
$ws.Cells.Item(43, 5).Value = "  +7.01%  "
$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).Value = "'18.39"
$ws.Cells.Item(44, 5).Value = "  +27.01%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).Value = "'0.998"
$ws.Cells.Item(45, 5).Value = "  -0.10%  "
$ws.Cells.Item(46, 4).Value = "'91.69"
$ws.Cells.Item(46, 5).Value = "  -1.48%  "
$ws.Cells.Item(47, 4).Value = "'1.79"
$ws.Cells.Item(47, 5).Value = "  -1.09%  "
$ws.Cells.Item(48, 4).Value = "'108.73"
$ws.Cells.Item(48, 5).Value = "  +11.23%  "
$ws.Cells.Item(49, 4).Value = "'9.03"
$ws.Cells.Item(49, 5).Value = "  +9.32%  "
$ws.Cells.Item(50, 5).Value = "  +7.45%  "
$ws.Cells.Item(51, 4).Value = "2.835.37"
$ws.Cells.Item(51, 5).Value = "  +9.57%  "
